$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = "Error from literal:"
$ws.Range("C21").Value = "#VALUE!"

$ws.Range("B22").Value = "Error from evaluation:"
$ws.Range("C22").Formula = "=1/0"
